# Stundenliste-SebastianEderer.xlsx edit: "Integrated SASS in WebUI"
# Adds a new timesheet entry on row 18 (date 2021-01-25, 5 hours,
# "Development WebUI") which cascades the running-total formula in
# column C from 51 to 56 for every subsequent row, and moves the
# active selection to I13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 18 values.
$ws.Range("B18").Value = 5
$ws.Range("D18").Value = "Development WebUI"

# A18 needs the same date-number-format style as the cells above it
# (numFmtId 14), so copy the format from A17 and then set the value,
# instead of letting COM infer/create a brand new date-time style.
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A18").Value = 44221

# Move the active selection as recorded in the saved sheet view.
$ws.Range("I13").Select()
